# Update the "Chiffres COVID-19 Valais" tracking sheet with the latest
# day's figures and correct a handful of previously-entered raw inputs.
# All cumulative/derived columns (B, H, J, K) are formulas and recalculate
# automatically; only the raw input columns are written here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Selection / view change -------------------------------------------------
# Select the merged title row instead of the previous "A2" selection.
$ws.Range("A1:M1").Select()

# --- Corrections to previously entered "new cases" (column C) ---------------
$ws.Range("C244").Value = 883
$ws.Range("C275").Value = 140
$ws.Range("C315").Value = 168
$ws.Range("C380").Value = 77
$ws.Range("C615").Value = 50
$ws.Range("C623").Value = 64
$ws.Range("C641").Value = 154
$ws.Range("C644").Value = 324
$ws.Range("C646").Value = 341
$ws.Range("C648").Value = 141
$ws.Range("C650").Value = 474
$ws.Range("C651").Value = 379
$ws.Range("C655").Value = 245
$ws.Range("C657").Value = 535
$ws.Range("C658").Value = 430

# --- Row 659 (2021-12-15): revised deaths split ------------------------------
$ws.Range("C659").Value = 378
$ws.Range("L659").Value = 1

# --- Row 660 (2021-12-16): revised cases + an extra-hospital death ----------
$ws.Range("C660").Value = 384
$ws.Range("M660").Value = 1

# --- Rows 661-664: newly reported days ---------------------------------------
$ws.Range("C661").Value = 338
$ws.Range("E661").Value = 10
$ws.Range("F661").Value = 4
$ws.Range("G661").Value = 71
$ws.Range("L661").Value = 0
$ws.Range("M661").Value = 1

$ws.Range("C662").Value = 187
$ws.Range("E662").Value = 8
$ws.Range("F662").Value = 5
$ws.Range("G662").Value = 74
$ws.Range("L662").Value = 0
$ws.Range("M662").Value = 0

$ws.Range("C663").Value = 96
$ws.Range("E663").Value = 8
$ws.Range("F663").Value = 4
$ws.Range("G663").Value = 77
$ws.Range("L663").Value = 0
$ws.Range("M663").Value = 0

$ws.Range("C664").Value = 17
$ws.Range("E664").Value = 7
$ws.Range("F664").Value = 4
$ws.Range("G664").Value = 78
$ws.Range("L664").Value = 0
$ws.Range("M664").Value = 0
